$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.1296102336659
$ws.Range("C2").Value = 12.51061022566877
$ws.Range("E2").Value = 15.39268966090119
$ws.Range("F2").Value = 45.50279647625896
$ws.Range("G2").Value = 3.672419454595956
$ws.Range("I2").Value = 27.48889242981799
$ws.Range("J2").Value = 9.026419800367718
$ws.Range("M2").Value = 18.91083277741815
$ws.Range("B3").Value = 14.66830947948019
$ws.Range("C3").Value = 11.98712236882445
$ws.Range("E3").Value = 15.41986848673057
$ws.Range("F3").Value = 45.21257085386257
$ws.Range("G3").Value = 3.676596264746986
$ws.Range("I3").Value = 27.36398792096362
$ws.Range("J3").Value = 9.044798418533569
$ws.Range("M3").Value = 18.82739027006425
$ws.Range("B4").Value = 14.38280619574022
$ws.Range("C4").Value = 11.65810411037496
$ws.Range("E4").Value = 15.43814699484551
$ws.Range("F4").Value = 45.04541914936755
$ws.Range("G4").Value = 3.679290494736608
$ws.Range("I4").Value = 27.2935947018549
$ws.Range("J4").Value = 9.057070311397508
$ws.Range("M4").Value = 18.78083309253585
$ws.Range("B5").Value = 14.26610190406804
$ws.Range("C5").Value = 11.52234809200102
$ws.Range("E5").Value = 15.44599704048196
$ws.Range("F5").Value = 44.98012685559629
$ws.Range("G5").Value = 3.680421151240829
$ws.Range("I5").Value = 27.26650108907947
$ws.Range("J5").Value = 9.062318966659099
$ws.Range("M5").Value = 18.76305276533897
$ws.Range("B6").Value = 14.24670790055
$ws.Range("C6").Value = 11.49971180779482
$ws.Range("E6").Value = 15.44732482070983
$ws.Range("F6").Value = 44.96945686697436
$ws.Range("G6").Value = 3.680610876988667
$ws.Range("I6").Value = 27.26209857648393
$ws.Range("J6").Value = 9.063205450820799
$ws.Range("M6").Value = 18.76017279379496
$ws.Range("B7").Value = 14.38123345019184
$ws.Range("C7").Value = 11.65627972815113
$ws.Range("E7").Value = 15.43825123607785
$ws.Range("F7").Value = 45.04452710717584
$ws.Range("G7").Value = 3.679305610444267
$ws.Range("I7").Value = 27.29322285215597
$ws.Range("J7").Value = 9.057140094133228
$ws.Range("M7").Value = 18.78058845407551
$ws.Range("B8").Value = 14.9711452629486
$ws.Range("C8").Value = 12.33181866062376
$ws.Range("E8").Value = 15.40173165105663
$ws.Range("F8").Value = 45.4004623667247
$ws.Range("G8").Value = 3.67383279037965
$ws.Range("I8").Value = 27.44452378823417
$ws.Range("J8").Value = 9.032551513050311
$ws.Range("M8").Value = 18.88110123204518
$ws.Range("B9").Value = 16.10145553282693
$ws.Range("C9").Value = 13.58742006627016
$ws.Range("E9").Value = 15.34267475066007
$ws.Range("F9").Value = 46.1838890171344
$ws.Range("G9").Value = 3.664123215144903
$ws.Range("I9").Value = 27.79072959967322
$ws.Range("J9").Value = 8.992190787705342
$ws.Range("M9").Value = 19.11453757372785
$ws.Range("B10").Value = 16.90528232782922
$ws.Range("C10").Value = 14.45759293425087
$ws.Range("E10").Value = 15.30685728316319
$ws.Range("F10").Value = 46.80827174695793
$ws.Range("G10").Value = 3.657604419982888
$ws.Range("I10").Value = 28.07448810345904
$ws.Range("J10").Value = 8.967359926920917
$ws.Range("M10").Value = 19.3070520199498
$ws.Range("B11").Value = 17.26327036128948
$ws.Range("C11").Value = 14.84039081089569
$ws.Range("E11").Value = 15.2921894037814
$ws.Range("F11").Value = 47.10213036550553
$ws.Range("G11").Value = 3.654770516853836
$ws.Range("I11").Value = 28.20974943877492
$ws.Range("J11").Value = 8.957117794722659
$ws.Range("M11").Value = 19.39892306018782
$ws.Range("B12").Value = 17.3975877047657
$ws.Range("C12").Value = 14.98335401764389
$ws.Range("E12").Value = 15.28686738808811
$ws.Range("F12").Value = 47.21474813226167
$ws.Range("G12").Value = 3.653716163213058
$ws.Range("I12").Value = 28.2618365088951
$ws.Range("J12").Value = 8.95339145866763
$ws.Range("M12").Value = 19.43430528938819
$ws.Range("B13").Value = 17.3687175247961
$ws.Range("C13").Value = 14.95265471161505
$ws.Range("E13").Value = 15.28800326556906
$ws.Range("F13").Value = 47.19043547805441
$ws.Range("G13").Value = 3.653942403775625
$ws.Range("I13").Value = 28.25058045534619
$ws.Range("J13").Value = 8.95418721526914
$ws.Range("M13").Value = 19.42665914311951
$ws.Range("B14").Value = 17.2743464376316
$ws.Range("C14").Value = 14.85219305689839
$ws.Range("E14").Value = 15.29174690872837
$ws.Range("F14").Value = 47.11136907056639
$ws.Range("G14").Value = 3.654683398829548
$ws.Range("I14").Value = 28.21401743875028
$ws.Range("J14").Value = 8.956808174092746
$ws.Range("M14").Value = 19.40182224641753
$ws.Range("B15").Value = 17.21637535114519
$ws.Range("C15").Value = 14.79039450945511
$ws.Range("E15").Value = 15.2940702217975
$ws.Range("F15").Value = 47.06311089710734
$ws.Range("G15").Value = 3.655139722073367
$ws.Range("I15").Value = 28.19173370229262
$ws.Range("J15").Value = 8.958433418786225
$ws.Range("M15").Value = 19.3866853528667
$ws.Range("B16").Value = 16.88172010499234
$ws.Range("C16").Value = 14.43230316788135
$ws.Range("E16").Value = 15.30784845828671
$ws.Range("F16").Value = 46.78925892198313
$ws.Range("G16").Value = 3.657792257825503
$ws.Range("I16").Value = 28.06577095183214
$ws.Range("J16").Value = 8.968050528916999
$ws.Range("M16").Value = 19.3011323500558
$ws.Range("B17").Value = 16.67434878771885
$ws.Range("C17").Value = 14.20919441248115
$ws.Range("E17").Value = 15.31671635044586
$ws.Range("F17").Value = 46.62372489335846
$ws.Range("G17").Value = 3.659453098013406
$ws.Range("I17").Value = 27.99006431663538
$ws.Range("J17").Value = 8.974220606708107
$ws.Range("M17").Value = 19.24973102204295
$ws.Range("B18").Value = 16.55435921062008
$ws.Range("C18").Value = 14.07964586511583
$ws.Range("E18").Value = 15.3219700248126
$ws.Range("F18").Value = 46.52944331811538
$ws.Range("G18").Value = 3.660420757289546
$ws.Range("I18").Value = 27.94710281829235
$ws.Range("J18").Value = 8.97786859912466
$ws.Range("M18").Value = 19.22057219308313
$ws.Range("B19").Value = 16.5136145202082
$ws.Range("C19").Value = 14.03557689830606
$ws.Range("E19").Value = 15.32377516536338
$ws.Range("F19").Value = 46.49768295929788
$ws.Range("G19").Value = 3.660750521688484
$ws.Range("I19").Value = 27.93265755094569
$ws.Range("J19").Value = 8.979120752422302
$ws.Range("M19").Value = 19.21076996506463
$ws.Range("B20").Value = 16.69649889396416
$ws.Range("C20").Value = 14.23307213098726
$ws.Range("E20").Value = 15.3157565131903
$ws.Range("F20").Value = 46.64125062039796
$ws.Range("G20").Value = 3.659275017578358
$ws.Range("I20").Value = 27.99806323583035
$ws.Range("J20").Value = 8.973553527536099
$ws.Range("M20").Value = 19.25516095445309
$ws.Range("B21").Value = 17.30210035418184
$ws.Range("C21").Value = 14.8817560510654
$ws.Range("E21").Value = 15.29064101388013
$ws.Range("F21").Value = 47.13455702092645
$ws.Range("G21").Value = 3.654465241861787
$ws.Range("I21").Value = 28.22473354920114
$ws.Range("J21").Value = 8.956034200892638
$ws.Range("M21").Value = 19.40910156456762
$ws.Range("B22").Value = 17.69058393044493
$ws.Range("C22").Value = 15.29404280804348
$ws.Range("E22").Value = 15.2755805067736
$ws.Range("F22").Value = 47.4647354929114
$ws.Range("G22").Value = 3.651431202642397
$ws.Range("I22").Value = 28.37791509003253
$ws.Range("J22").Value = 8.945471378417034
$ws.Range("M22").Value = 19.5131528552337
$ws.Range("B23").Value = 17.4839549004983
$ws.Range("C23").Value = 15.07509956268644
$ws.Range("E23").Value = 15.28349515198393
$ws.Range("F23").Value = 47.28782643688461
$ws.Range("G23").Value = 3.653040556706963
$ws.Range("I23").Value = 28.29570582601387
$ws.Range("J23").Value = 8.951027576837532
$ws.Range("M23").Value = 19.45731227120647
$ws.Range("B24").Value = 16.68648721547613
$ws.Range("C24").Value = 14.22228099334733
$ws.Range("E24").Value = 15.31618997145682
$ws.Range("F24").Value = 46.63332447358253
$ws.Range("G24").Value = 3.659355487801391
$ws.Range("I24").Value = 27.99444516841359
$ws.Range("J24").Value = 8.973854800318861
$ws.Range("M24").Value = 19.25270485764245
$ws.Range("B25").Value = 15.79965049599985
$ws.Range("C25").Value = 13.25628672353208
$ws.Range("E25").Value = 15.35731550193268
$ws.Range("F25").Value = 45.96312672849504
$ws.Range("G25").Value = 3.666641323154376
$ws.Range("I25").Value = 27.69184324223808
$ws.Range("J25").Value = 9.002264943408175
$ws.Range("M25").Value = 19.04760949741692
